# Apply updated loading_percent results for the 380 kV case (rows 2-25, columns B-O)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "B2" = 8.706829043602321
    "C2" = 5.786482077707481
    "D2" = 5.972618362059277
    "E2" = 16.50508050714844
    "G2" = 24.96442592711288
    "H2" = 13.27028801887965
    "I2" = 18.79813087760017
    "K2" = 8.459906954226811
    "N2" = 17.31343588228232
    "O2" = 19.70215311645135
    "B3" = 8.35495519919469
    "C3" = 5.551226904570221
    "D3" = 5.851287856702156
    "E3" = 15.56961735666164
    "G3" = 24.94012192324781
    "H3" = 13.31010655170421
    "I3" = 18.87449877332737
    "K3" = 8.176206086240951
    "N3" = 17.36943826834594
    "O3" = 19.7520826941456
    "B4" = 8.132324813206095
    "C4" = 5.4004323931131
    "D4" = 5.777242512512783
    "E4" = 14.97024299731406
    "G4" = 24.9343423236375
    "H4" = 13.33681765142158
    "I4" = 18.92519307239155
    "K4" = 8.018053614467387
    "N4" = 17.40541859906455
    "O4" = 19.78740714480531
    "B5" = 8.040084310019422
    "C5" = 5.337437643761288
    "D5" = 5.747229564753606
    "E5" = 14.71997686137593
    "G5" = 24.93428411401652
    "H5" = 13.34827098492219
    "I5" = 18.94680641705234
    "K5" = 7.952994473427749
    "N5" = 17.4204831056166
    "O5" = 19.80297234845015
    "B6" = 8.024680565922436
    "C6" = 5.326885752302799
    "D6" = 5.742257067926505
    "E6" = 14.67806591014319
    "G6" = 24.9344130723452
    "H6" = 13.35020711259264
    "I6" = 18.95045293094222
    "K6" = 7.942157605227822
    "N6" = 17.42300888575959
    "O6" = 19.80562751009105
    "B7" = 8.131086773628869
    "C7" = 5.399589005614864
    "D7" = 5.776837031832853
    "E7" = 14.96689178545909
    "G7" = 24.93433224274304
    "H7" = 13.33696981463068
    "I7" = 18.92548069281999
    "K7" = 8.017178540582609
    "N7" = 17.40562013401573
    "O7" = 19.78761232927212
    "B8" = 8.586945780794888
    "C8" = 5.706711937203442
    "D8" = 5.93071802524158
    "E8" = 16.18786669062461
    "G8" = 24.95414774158237
    "H8" = 13.28354757293729
    "I8" = 18.82367204449866
    "K8" = 8.342599826360436
    "N8" = 17.33241528991127
    "O8" = 19.71839807359844
    "B9" = 9.42344404410616
    "C9" = 6.25671567658154
    "D9" = 6.234023935146902
    "E9" = 18.46251569797813
    "G9" = 25.06552985374351
    "H9" = 13.19676164188146
    "I9" = 18.65427962558202
    "K9" = 9.15348828553679
    "N9" = 17.2014558862951
    "O9" = 19.61985028312212
    "B10" = 9.997004379490994
    "C10" = 6.626823079600021
    "D10" = 6.455163977654201
    "E10" = 20.10046457254134
    "G10" = 25.19137517881335
    "H10" = 13.14399058153259
    "I10" = 18.54836699964218
    "K10" = 9.698677715222638
    "N10" = 17.11283632084295
    "O10" = 19.57030136615639
    "B11" = 10.24802187577084
    "C11" = 6.787474359084514
    "D11" = 6.554873343825821
    "E11" = 20.80331738959663
    "G11" = 25.25808474121488
    "H11" = 13.12237690059124
    "I11" = 18.50423137172772
    "K11" = 9.935154802710043
    "N11" = 17.07415296065296
    "O11" = 19.55275768468358
    "B12" = 10.3415839486686
    "C12" = 6.847177014331014
    "D12" = 6.592462205890083
    "E12" = 21.0634285100052
    "G12" = 25.28469410457518
    "H12" = 13.11453689758375
    "I12" = 18.48810177683492
    "K12" = 10.02301175003369
    "N12" = 17.05973769439581
    "O12" = 19.54683518969851
    "B13" = 10.32150107554754
    "C13" = 6.83436968059045
    "D13" = 6.584374914126327
    "E13" = 21.00767721580928
    "G13" = 25.27890358361127
    "H13" = 13.11621004503963
    "I13" = 18.49154958545633
    "K13" = 10.00416590169445
    "N13" = 17.06283192242942
    "O13" = 19.54807860636315
    "B14" = 10.25574956429901
    "C14" = 6.792408950218955
    "D14" = 6.55796941460934
    "E14" = 20.82483783200944
    "G14" = 25.26024697759408
    "H14" = 13.12172498834681
    "I14" = 18.50289267232268
    "K14" = 9.942416924718609
    "N14" = 17.07296233910916
    "O14" = 19.55225597631393
    "B15" = 10.21527856711738
    "C15" = 6.766558675828972
    "D15" = 6.541772105697325
    "E15" = 20.71205729123444
    "G15" = 25.24899438280861
    "H15" = 13.12514794991481
    "I15" = 18.50991670343474
    "K15" = 9.904372673643321
    "N15" = 17.0791978588536
    "O15" = 19.55490868791685
    "B16" = 9.980394450358965
    "C16" = 6.616166852551499
    "D16" = 6.448626286727871
    "E16" = 20.05368394371553
    "G16" = 25.18720492230727
    "H16" = 13.14545127773775
    "I16" = 18.55133292494926
    "K16" = 9.682988332590604
    "N16" = 17.11539707507957
    "O16" = 19.57154864269688
    "B17" = 9.833712827389302
    "C17" = 6.521913057345463
    "D17" = 6.39122677429404
    "E17" = 19.63899319419209
    "G17" = 25.15171429078761
    "H17" = 13.15851985945462
    "I17" = 18.57777776906357
    "K17" = 9.544196884115618
    "N17" = 17.13802083037933
    "O17" = 19.58303824239499
    "B18" = 9.748417256931031
    "C18" = 6.466975973257634
    "D18" = 6.358130793200916
    "E18" = 19.39649454757093
    "G18" = 25.13219254385127
    "H18" = 13.16626168425881
    "I18" = 18.59336871612144
    "K18" = 9.463283852148056
    "N18" = 17.15118691775646
    "O18" = 19.59011691095038
    "B19" = 9.719380511771435
    "C19" = 6.448251515428809
    "D19" = 6.346912392647604
    "E19" = 19.31370418970688
    "G19" = 25.12573628090768
    "H19" = 13.16892157819671
    "I19" = 18.59871284554769
    "K19" = 9.435703048076027
    "N19" = 17.15567112732322
    "O19" = 19.5925942965654
    "B20" = 9.849423935245349
    "C20" = 6.532021730403306
    "D20" = 6.39734577039536
    "E20" = 19.68354935756973
    "G20" = 25.15540014332482
    "H20" = 13.15710538257473
    "I20" = 18.57492326792913
    "K20" = 9.559083906493452
    "N20" = 17.13559661433789
    "O20" = 19.58176647692808
    "B21" = 10.27510339545276
    "C21" = 6.804764737580812
    "D21" = 6.565730250250466
    "E21" = 20.87870596368125
    "G21" = 25.2656904074719
    "H21" = 13.1200957576699
    "I21" = 18.49954507772075
    "K21" = 9.960600250123068
    "N21" = 17.06998046787023
    "O21" = 19.55100939807952
    "B22" = 10.54458292162685
    "C22" = 6.976408710194169
    "D22" = 6.674777323039194
    "E22" = 21.62459592327487
    "G22" = 25.34562036413751
    "H22" = 13.09791687200166
    "I22" = 18.45368373203003
    "K22" = 10.21314156191944
    "N22" = 17.02845575755783
    "O22" = 19.53511083574695
    "B23" = 10.40157553188672
    "C23" = 6.885410632785498
    "D23" = 6.616681497845119
    "E23" = 21.22971174842338
    "G23" = 25.30224692880466
    "H23" = 13.10957012643675
    "I23" = 18.47784878688185
    "K23" = 10.07926853224856
    "N23" = 17.05049427662536
    "O23" = 19.54321089861266
    "B24" = 9.842323947962369
    "C24" = 6.527453931525259
    "D24" = 6.394579668870202
    "E24" = 19.66341825048681
    "G24" = 25.15373102011243
    "H24" = 13.15774415575333
    "I24" = 18.57621258110526
    "K24" = 9.552356969627521
    "N24" = 17.1366921058098
    "O24" = 19.58233996841506
    "B25" = 9.203970759141992
    "C25" = 6.113760550023247
    "D25" = 6.15208527461316
    "E25" = 17.82163492391986
    "G25" = 25.02764401372961
    "H25" = 13.21831177179093
    "I25" = 18.69685596433876
    "K25" = 8.942747372565419
    "N25" = 17.23554399989822
    "O25" = 19.64250917396901
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
